$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = "10 sets count"

$ws.Range("M2:M39").Formula = "=A2*10"

$ws.Range("N2").Select() | Out-Null
